# Commit: Bazar -213(oil300,lau,begun,alu,piyaj,motaDal,chotoDal), Taher+250, Forhad+250, Himel+200
#
# Sheet layout reminder:
#   Row 25 = Himel, Row 27 = Taher, Row 28 = Forhad  -> column O holds each
#            person's "deposit" entry for this bazar event.
#   Row 42 = per-bazar item/notes row, Row 43 = the bazar cost entered by the
#            person named in row 42.  Row 46 totals (dep/cost) recalc off of
#            these automatically, as do every other downstream formula
#            (AG/AH/AI/AJ columns), so only the raw inputs need to change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New deposits recorded against this bazar run
$ws.Range("O25").Value = 200   # Himel +200
$ws.Range("O27").Value = 250   # Taher +250
$ws.Range("O28").Value = 250   # Forhad +250

# Nayem did this bazar run (oil, lau, begun, alu, piyaj, motaDal, chotoDal)
$ws.Range("O42").Value = "Nayem"
$ws.Range("O43").Value = 213   # Bazar cost -213

# Leave the cursor where the author left it when saving
$ws.Range("O29").Select()

$wb.Save()
